$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 2 ("1 | asd | asd | asd") -> last cell "asd" becomes "sdasd"
$table.Cell(2, 4).Range.Text = "sdasd"

# Row 3 ("2 | hjk | hjk | hjk") -> "cotizacion2" / "800" / "600"
$table.Cell(3, 2).Range.Text = "cotizacion2"
$table.Cell(3, 3).Range.Text = "800"
$table.Cell(3, 4).Range.Text = "600"

# Row 4 ("3 | klñ | klñ | klñ") is removed entirely
$table.Rows.Item(4).Delete()
